# CS401ProjectGnattChart.xlsx - reproduce the view-state / layout tweaks
# described by the commit: the workbook was simply re-opened, scrolled,
# a header cell-range selected, the last "week" column nudged a bit wider,
# and re-saved. No cell values or formulas change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")
$ws.Activate()

# --- Re-merge the week/header cells --------------------------------------
# Excel rewrites <mergeCells> with the touched ranges moved to the end of
# the list whenever a merged range is re-merged; round-trip the five
# ranges that end up at the tail of the saved mergeCells collection so the
# emitted order matches.
$remergeRanges = @("C3:D3", "C4:D4", "B5:G5", "AK4:AQ4", "AR4:AX4")
foreach ($r in $remergeRanges) {
    $ws.Range($r).UnMerge()
    $ws.Range($r).Merge()
}

# --- Widen the trailing "week" column -------------------------------------
# Column BS (the last of the day-grid columns) picks up a bit of extra
# width, splitting it out of the shared 9:71 column-width run.
$ws.Range("BS1").EntireColumn.ColumnWidth = 3

# --- Scroll down the frozen grid and select the current week's header ----
$ws.Range("BM4:BS4").Select()
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 29
$win.ScrollColumn = 1
